$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update SmokeTest (column D) value from "Yes" to "No" for rows 336-370
$ws.Range("D336:D370").Value = "No"

# Update the view: scroll position and selection to match the saved state
$ws.Range("C336").Select()
$excel.ActiveWindow.ScrollRow = 324
